$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Rename Sheet2 -> Data, make it the active/selected sheet ---
$ws2.Name = "Data"
$ws2.Activate()

# --- Column A: labels (written so shared-string table fills up in the
#     same order the original workbook used) ---
$ws2.Range("A1").Value = "Int"
$ws2.Range("A2").Value = "Double"
$ws2.Range("A3").Value = "Double Formula"
$ws2.Range("A4").Value = "String"
$ws2.Range("B4").Value = "Test"
$ws2.Range("A5").Value = "Percent"
$ws2.Range("A6").Value = "String Formular"
$ws2.Range("A7").Value = "Date"
$ws2.Range("A9").Value = "Date Formular"
$ws2.Range("A8").Value = "Date Func"

# --- Column B: sample values exercising each data/formula type ---

# Int
$ws2.Range("B1").Value = 1

# Double
$ws2.Range("B2").Value = 0.25

# Double Formula
$ws2.Range("B3").Formula = "=6/10"

# Percent (set format before the value so the engine doesn't fabricate
# a throwaway custom numFmt first)
$ws2.Range("B5").NumberFormat = "0%"
$ws2.Range("B5").Value = 0.1

# String Formular
$ws2.Range("B6").Formula = '="A" & "B"'

# Date (literal serial value, d-mmm display)
$ws2.Range("B7").NumberFormat = "d-mmm"
$ws2.Range("B7").Value = 40939

# Date Func / Date Formular: apply the format first, then the formula,
# then clone the resulting style onto the dependent cell via
# copy/paste-special so both cells share a single cellXf (mirrors how
# Excel itself dedupes identical formats) instead of minting a second,
# redundant one.
$ws2.Range("B8").NumberFormat = "mm-dd-yy"
$ws2.Range("B8").Formula = "=TODAY()"

$null = $ws2.Range("B8").Copy()
$null = $ws2.Range("B9").PasteSpecial(-4122)
$ws2.Range("B9").Formula = "=B8+1"

# --- Column widths for the new sheet ---
$ws2.Columns.Item(1).ColumnWidth = 14.5
$ws2.Columns.Item(2).ColumnWidth = 13

# --- Selection matches the captured state ---
$null = $ws2.Range("G19").Select()

Write-Output "done"
